$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
# "Volume 31   Number  18" -> "...19"
$ws.Range("A8").Value = "Volume 31   Number  19"
# "Report Covering the Week  4/29/2024  Through  5/5/2024" -> new dates
$ws.Range("C9").Value = "Report Covering the Week  5/6/2024  Through  5/12/2024"

# --- Column width change (col H, index 8) ---
$ws.Columns.Item(8).ColumnWidth = 6.168446

# --- Row 14 (Murder) ---
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = "***.*"

# --- Row 15 (Rape) ---
$ws.Range("F15").Value = 2
$ws.Range("I15").Value = 5
$ws.Range("K15").Value = 150
$ws.Range("L15").Value = 150
$ws.Range("M15").Value = 66.666666666666
$ws.Range("N15").Value = -58.333333333333

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -25
$ws.Range("I16").Value = 68
$ws.Range("J16").Value = 62
$ws.Range("K16").Value = 9.677419354838
$ws.Range("L16").Value = 11.475409836065
$ws.Range("M16").Value = -4.225352112676
$ws.Range("N16").Value = -75.971731448763

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 25
$ws.Range("H17").Value = 66.666666666666
$ws.Range("I17").Value = 95
$ws.Range("J17").Value = 94
$ws.Range("K17").Value = 1.063829787234
$ws.Range("L17").Value = 14.457831325301
$ws.Range("M17").Value = 41.791044776119
$ws.Range("N17").Value = -57.399103139013

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -100
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -57.142857142857
$ws.Range("J18").Value = 45
$ws.Range("K18").Value = -31.111111111111
$ws.Range("L18").Value = -55.072463768115
$ws.Range("M18").Value = -31.111111111111
$ws.Range("N18").Value = -90.220820189274

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 25
$ws.Range("H19").Value = 76
$ws.Range("I19").Value = 186
$ws.Range("J19").Value = 130
$ws.Range("K19").Value = 43.076923076923
$ws.Range("L19").Value = 47.619047619047
$ws.Range("M19").Value = 86
$ws.Range("N19").Value = 25.675675675675

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = "***.*"
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 19
$ws.Range("K20").Value = -5
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 216.666666666667
$ws.Range("N20").Value = -66.666666666666

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -4
$ws.Range("F21").Value = 91
$ws.Range("G21").Value = 67
$ws.Range("H21").Value = 35.820895522388
$ws.Range("I21").Value = 406
$ws.Range("J21").Value = 356
$ws.Range("K21").Value = 14.044943820224
$ws.Range("L21").Value = 12.465373961218
$ws.Range("M21").Value = 38.095238095238
$ws.Range("N21").Value = -61.406844106463

# --- Row 22 (Transit) ---
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 7
$ws.Range("J22").Value = 6
$ws.Range("K22").Value = 16.666666666666
$ws.Range("L22").Value = -53.333333333333
$ws.Range("M22").Value = -12.5

# --- Row 23 (Housing) ---
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("F23").Value = 9
$ws.Range("H23").Value = 80
$ws.Range("I23").Value = 35
$ws.Range("J23").Value = 32
$ws.Range("K23").Value = 9.375
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 218.181818181818

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = -38.235294117647
$ws.Range("F24").Value = 125
$ws.Range("G24").Value = 105
$ws.Range("H24").Value = 19.047619047619
$ws.Range("I24").Value = 455
$ws.Range("J24").Value = 441
$ws.Range("K24").Value = 3.174603174603
$ws.Range("L24").Value = 15.776081424936
$ws.Range("M24").Value = 34.218289085545

# --- Row 25 (Retail Theft) ---
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 19
$ws.Range("E25").Value = -47.368421052631
$ws.Range("F25").Value = 80
$ws.Range("G25").Value = 61
$ws.Range("H25").Value = 31.147540983606
$ws.Range("I25").Value = 286
$ws.Range("J25").Value = 270
$ws.Range("K25").Value = 5.925925925925
$ws.Range("L25").Value = 40.196078431372

# --- Row 26 (Misd. Assault) ---
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = -68.75
$ws.Range("F26").Value = 29
$ws.Range("G26").Value = 46
$ws.Range("H26").Value = -36.956521739130
$ws.Range("I26").Value = 133
$ws.Range("J26").Value = 153
$ws.Range("K26").Value = -13.071895424836
$ws.Range("L26").Value = -16.352201257861
$ws.Range("M26").Value = -23.121387283237

# --- Row 27 (UCR Rape*) ---
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 7
$ws.Range("K27").Value = 40
$ws.Range("L27").Value = 75

# --- Row 28 (Other Sex Crimes) ---
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -25
$ws.Range("J28").Value = 15
$ws.Range("K28").Value = 6.666666666666
$ws.Range("L28").Value = -27.272727272727

# --- Row 29 (Shooting Vic.) ---
$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = "***.*"
$ws.Range("I29").Value = 7
$ws.Range("K29").Value = 16.666666666666
$ws.Range("L29").Value = 16.666666666666
$ws.Range("M29").Value = -12.5
$ws.Range("N29").Value = -80.555555555555

# --- Row 30 (Shooting Inc.) ---
$ws.Range("C30").Value = 1
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = "***.*"
$ws.Range("I30").Value = 6
$ws.Range("K30").Value = 20
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -25
$ws.Range("N30").Value = -80.645161290322
